$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.241.95'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.775.91'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '633.55'
$ws.Range('E5').Value = '  +4.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.68'
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.773.54'
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.460'
$ws.Range('E11').Value = '  +2.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.76'
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.99'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.409.03'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.779.11'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.206.11'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.03'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '463.93'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.55'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.78'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000145'
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.924.48'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.31'
$ws.Range('E31').Value = '  +4.53%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.69'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.52'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.167'
$ws.Range('E36').Value = '  +14.16%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.727.22'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.98'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.31'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '158.18'
$ws.Range('E45').Value = '  +3.35%  '
$ws.Range('E46').Value = '  +5.93%  '
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '43.29'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.78'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('E51').Value = '  +0.46%  '
